$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column C
$ws.Range("C1").Value = "z"

# C2: direct (non-shared) formula, mirrors B2's pattern
$ws.Range("C2").Formula = "=A2^3-5*A2^2"

# C3:C10: shared formula, mirrors B3:B10's pattern
$ws.Range("C3:C10").Formula = "=A3^3-5*A3^2"

# Column C width ~10 (matches the author's bestFit column sizing)
$ws.Columns("C").ColumnWidth = 9.166666666666666

# Active cell moves to C2, as in the edited file
$ws.Range("C2").Select()
